$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-19 Sunday" "2024-05-20 Monday"

Replace-Text "15×38=" "34×38="
Replace-Text "75×37=" "39×49="
Replace-Text "66×32=" "74×88="
Replace-Text "21×95=" "84×42="
Replace-Text "14×61=" "47×94="

Replace-Text "41×42=" "78×53="
Replace-Text "62×29=" "60×50="
Replace-Text "62×60=" "19×39="
Replace-Text "83×84=" "38×61="
Replace-Text "74×60=" "65×72="

Replace-Text "57×19=" "82×84="
Replace-Text "61×52=" "38×83="
Replace-Text "69×20=" "29×38="
Replace-Text "12×98=" "95×91="
Replace-Text "55×45=" "74×35="

Replace-Text "16×44=" "15×96="
Replace-Text "51×15=" "96×86="
Replace-Text "41×60=" "72×28="
Replace-Text "86×87=" "48×98="
Replace-Text "80×29=" "55×78="

Replace-Text "83×28=" "31×65="
Replace-Text "56×48=" "65×60="
Replace-Text "70×99=" "80×75="
Replace-Text "83×23=" "84×73="
Replace-Text "71×22=" "22×42="
